# "Updates for Night of the Big Wind"
# Adds a new location row (Neath) to the rough-locations sheet and moves
# the selection the way it was left after the edit (scrolled down near
# the bottom of the list with E65 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended after the existing last row (row 59 -> 60):
# LOCATION, LAT, LON, NDAY, TYPE (1=exist,2=newspaper,3=diary)
$ws.Range("A60").Value = "Neath"
$ws.Range("B60").Value = 51.66
$ws.Range("C60").Value = -3.81
$ws.Range("D60").Value = 3
$ws.Range("E60").Value = 2

# Scroll the view down toward the bottom of the sheet and leave the
# selection where the author left it after adding the new row.
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E65").Select()
